$wb = $excel.ActiveWorkbook

# --- Update selection on "Site Map Done" (sheet2) to E7 ---
$wsSiteMap = $wb.Worksheets.Item("Site Map Done")
$wsSiteMap.Activate()
$wsSiteMap.Range("E7").Select()

# --- Add a new worksheet "Done Complete" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDone = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsDone.Name = "Done Complete"

# --- Populate the new sheet with its values ---
$values = @(
    404,
    503,
    "about",
    "about-solar",
    "blog",
    "blog-post",
    "brochures-technical-info",
    "careers",
    "case",
    "clients",
    "coming-soon",
    "commercial-government",
    "compare-plans-offers"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $wsDone.Cells.Item($i + 1, 1).Value = $values[$i]
}

# --- Make the new sheet the active sheet/tab and set its cursor position ---
$wsDone.Activate()
$wsDone.Range("A14").Select()
